$d = $word.ActiveDocument
$d.Content.Find.Execute("48÷7=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "91÷2=45, 1", 2) | Out-Null
$d.Content.Find.Execute("72÷9=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "20÷3=6, 2", 2) | Out-Null
$d.Content.Find.Execute("56÷3=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=42, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷2=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷7=7, 5", 2) | Out-Null
$d.Content.Find.Execute("77÷2=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=3, 1", 2) | Out-Null
$d.Content.Find.Execute("55÷2=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "80÷9=8, 8", 2) | Out-Null
$d.Content.Find.Execute("11÷9=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "19÷3=6, 1", 2) | Out-Null
$d.Content.Find.Execute("74÷4=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=3, 5", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "39÷6=6, 3", 2) | Out-Null
$d.Content.Find.Execute("14÷5=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "40÷9=4, 4", 2) | Out-Null
$d.Content.Find.Execute("70÷7=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "59÷8=7, 3", 2) | Out-Null
$d.Content.Find.Execute("26÷8=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=3, 1", 2) | Out-Null
$d.Content.Find.Execute("44÷2=22, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷7=12, 6", 2) | Out-Null
$d.Content.Find.Execute("25÷6=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "31÷4=7, 3", 2) | Out-Null
$d.Content.Find.Execute("71÷5=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "33÷3=11, 0", 2) | Out-Null
$d.Content.Find.Execute("15÷5=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=17, 2", 2) | Out-Null
$d.Content.Find.Execute("29÷5=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "81÷7=11, 4", 2) | Out-Null
$d.Content.Find.Execute("61÷4=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "79÷2=39, 1", 2) | Out-Null
$d.Content.Find.Execute("57÷4=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "18÷3=6, 0", 2) | Out-Null
$d.Content.Find.Execute("92÷6=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "60÷7=8, 4", 2) | Out-Null
$d.Content.Find.Execute("62÷6=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "74÷5=14, 4", 2) | Out-Null
$d.Content.Find.Execute("11÷3=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "62÷5=12, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "44÷6=7, 2", 2) | Out-Null
$d.Content.Find.Execute("19÷8=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "36÷3=12, 0", 2) | Out-Null
